$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.216.68'
$ws.Range('E2').Value = '  -0.12%  '
$ws.Range('D3').Value = '3.108.23'
$ws.Range('E3').Value = '  -0.16%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '574.82'
$ws.Range('E5').Value = '  -0.84%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '177.89'
$ws.Range('E6').Value = '  +2.78%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D8').Value = '3.109.65'
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.513'
$ws.Range('E9').Value = '  -1.31%  '
$ws.Range('E10').Value = '  -2.11%  '
$ws.Range('E11').Value = '  -0.74%  '
$ws.Range('E12').Value = '  -1.86%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000240'
$ws.Range('E13').Value = '  -2.78%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.09'
$ws.Range('E14').Value = '  -1.58%  '
$ws.Range('E15').Value = '  -0.22%  '
$ws.Range('D16').Value = '3.623.18'
$ws.Range('E16').Value = '  -0.16%  '
$ws.Range('D17').Value = '67.057.51'
$ws.Range('E17').Value = '  -0.18%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.04'
$ws.Range('E18').Value = '  -0.76%  '
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').Value = '3.095.31'
$ws.Range('E19').Value = '  -0.49%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.67'
$ws.Range('E20').Value = '  +0.29%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '485.71'
$ws.Range('E21').Value = '  -1.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.80'
$ws.Range('E22').Value = '  -0.21%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.689'
$ws.Range('E23').Value = '  -1.59%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '83.46'
$ws.Range('E24').Value = '  -0.44%  '
$ws.Range('B25').Value = 'Fetch.AI'
$ws.Range('C25').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.27'
$ws.Range('E25').Value = '  -0.96%  '
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.61'
$ws.Range('E26').Value = '  -3.80%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.14'
$ws.Range('E27').Value = '  -4.04%  '
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('E29').Value = '  +0.65%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.31'
$ws.Range('E30').Value = '  -1.56%  '
$ws.Range('E31').Value = '  -2.66%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '28.11'
$ws.Range('E32').Value = '  -0.45%  '
$ws.Range('E33').Value = '  -1.62%  '
$ws.Range('D34').Value = '0.0₃0942'
$ws.Range('E34').Value = '  -0.17%  '
$ws.Range('E35').Value = '  -0.09%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '48.32'
$ws.Range('E36').Value = '  +3.27%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.60'
$ws.Range('E37').Value = '  -3.29%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.948'
$ws.Range('E38').Value = '  -2.30%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.314'
$ws.Range('E39').Value = '  +2.02%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '49.12'
$ws.Range('E40').Value = '  -1.81%  '
$ws.Range('E41').Value = '  -0.69%  '
$ws.Range('E42').Value = '  -0.70%  '
$ws.Range('E43').Value = '  -1.59%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.73'
$ws.Range('E44').Value = '  +5.90%  '
$ws.Range('D45').Value = '2.794.08'
$ws.Range('E45').Value = '  -0.23%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '372.25'
$ws.Range('E46').Value = '  -3.36%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0345'
$ws.Range('E47').Value = '  -0.87%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '135.38'
$ws.Range('E48').Value = '  +0.05%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '25.49'
$ws.Range('E50').Value = '  +1.96%  '
$ws.Range('E51').Value = '  +4.38%  '
